# Restore revision: change cell C10 on the "Rules" sheet from 18 to 1
# (DESIGN/rules/Sample Project/Main.xlsx, row 10 "From" condition value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
